# Edit script: 
#  1) Fix the mislabeled title on the "Hardware" architecture slide
#     ("Architecture - Software" -> "Architecture - Hardware").
#  2) Swap the presentation's applied colour theme with the colour
#     scheme that used to live in the secondary (notes-only) theme
#     part -- i.e. the deck switches from the "SlateVTI" green theme
#     to the standard "Tema di Office" colours.

$p = $ppt.ActivePresentation

# --- 1) Title text fix on slide 3 -----------------------------------------
$slide3 = $p.Slides.Item(3)
$titleShape = $slide3.Shapes.Item(4)
$titleShape.TextFrame.TextRange.Text = "Architecture - Hardware"

# --- 2) Theme colour swap ---------------------------------------------------
# New ("Tema di Office") colour values, in the standard a:clrScheme order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$newColors = @(
    0,         # dk1     000000
    16777215,  # lt1     FFFFFF
    6968388,   # dk2     44546A
    15132391,  # lt2     E7E6E6
    12874308,  # accent1 4472C4
    3243501,   # accent2 ED7D31
    10855845,  # accent3 A5A5A5
    49407,     # accent4 FFC000
    13998939,  # accent5 5B9BD5
    4697456,   # accent6 70AD47
    12673797,  # hlink   0563C1
    7491477    # folHlink 954F72
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = $newColors[$i - 1]
}
